# DMASE - Test Results.xlsx : convert h-vector formulation from polar to
# rectangular on the "2 Bus" sheet, per the commit:
#   "Converted h vector formulation from polar to rectangular"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2 Bus")
$ws.Activate()

# --- Row 2-4: fill in the previously-empty L column (rectangular h residual)
$ws.Range("L2").Value = [double]"-4.09372882814729E-6"
$ws.Range("L3").Value = [double]"-5.8479390156740401E-6"
$ws.Range("L4").Value = [double]"-4.1007288289485197E-6"

# --- Row 5 ("Ang 2" stays) -------------------------------------------------
$ws.Range("B5").Value = [double]"-3.8391445791322899E-2"
$ws.Range("C5:G5").ClearContents()
$ws.Range("L5").Value = [double]"-5.8419390149555601E-6"

# --- Row 6 ("Ang 3" -> "V 1") ----------------------------------------------
$ws.Range("A6").Value = "V 1"
$ws.Range("B6").Value = [double]"1.00000000217498"
$ws.Range("C6:G6").ClearContents()
$ws.Range("L6").Value = [double]"3.76791437384227E-6"

# --- Row 7 ("V 1" -> "V 2") -------------------------------------------------
$ws.Range("A7").Value = "V 2"
$ws.Range("B7").Value = [double]"0.96399318557074398"
$ws.Range("C7:G7").ClearContents()
$ws.Range("L7").Value = [double]"5.3231864123454196E-6"

# --- Row 8 (label & A:G values removed, but cells remain, blank) ----------
$ws.Range("A8:G8").ClearContents()
$ws.Range("L8").Value = [double]"1.4224587097100299E-7"

# --- Row 9 (A:G cells removed entirely) ------------------------------------
$ws.Range("A9:G9").Clear()
$ws.Range("L9").Value = [double]"5.0522178718370704E-7"

# --- Selection moves from J5 to D12 ----------------------------------------
$ws.Range("D12").Select()

# --- Page setup: orientation now explicit (portrait) -----------------------
$ws.PageSetup.Orientation = 1

$wb.Save()
